# Add payment 79174445 (Cash) 2025-08-18T08:41:43
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns: birthday_discount (H), points_redeemed (I).
# Copy the existing header style (bold, bordered, centered) from A1 onto
# the new header cells so the style table isn't needlessly duplicated.
$ws.Range("A1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)
$ws.Cells.Item(1, 8).Value = "birthday_discount"
$ws.Cells.Item(1, 9).Value = "points_redeemed"

# Row 13: the phone number had been captured as text; normalize it to a
# real number like the rest of the `phone` column.
$ws.Cells.Item(13, 1).Value = 79172233

# New payment record in row 14.
# Keep the phone number as text (it's stored as text for new rows),
# so force a text number format before assigning, then drop back to the
# default "Normal" style so no stray formatting is left on the cell.
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = "79174445"
$ws.Cells.Item(14, 1).Style = "Normal"
$ws.Cells.Item(14, 3).Value = "Cash"
$ws.Cells.Item(14, 4).Value = "2025-08-18T08:41:43"
$ws.Cells.Item(14, 5).Value = 20
$ws.Cells.Item(14, 7).Value = 20
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0

$wb.Save()
